$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '52.111.98'
$ws.Range("E2").Value = '  +0.69%  '

$ws.Range("D3").Value = '2.916.87'
$ws.Range("E3").Value = '  +4.28%  '

$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").Value = "'354.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.90%  '

$ws.Range("D6").Value = "'113.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.69%  '

$ws.Range("D7").Value = "'0.558"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.10%  '

$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("D9").Value = "'0.623"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.39%  '

$ws.Range("D10").Value = "'39.57"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.40%  '

$ws.Range("E11").Value = '  +5.03%  '

$ws.Range("E12").Value = '  +1.06%  '

$ws.Range("E13").Value = '  +0.95%  '

$ws.Range("E14").Value = '  +0.06%  '

$ws.Range("D15").Value = '3.375.44'
$ws.Range("E15").Value = '  +4.26%  '

$ws.Range("D16").Value = '2.909.91'
$ws.Range("E16").Value = '  +3.94%  '

$ws.Range("D17").Value = "'0.987"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.39%  '

$ws.Range("D18").Value = '52.161.13'
$ws.Range("E18").Value = '  +0.82%  '

$ws.Range("E19").Value = '  -1.28%  '

$ws.Range("D20").Value = "'7.59"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.30%  '

$ws.Range("D21").Value = "'14.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.08%  '

$ws.Range("D22").Value = '0.0₃0979'
$ws.Range("E22").Value = '  +0.77%  '

$ws.Range("D23").Value = "'71.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.18%  '

$ws.Range("D24").Value = "'270.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.47%  '

$ws.Range("E25").Value = '  +2.09%  '

$ws.Range("E26").Value = '  +11.75%  '

$ws.Range("D27").Value = "'26.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.95%  '

$ws.Range("E28").Value = '  +0.03%  '

$ws.Range("D29").Value = "'7.04"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +15.29%  '

$ws.Range("E30").Value = '  +2.30%  '

$ws.Range("E31").Value = '  +13.09%  '

$ws.Range("B32").Value = 'Toncoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D32").Value = "'2.27"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.06%  '

$ws.Range("B33").Value = 'InjectiveProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D33").Value = "'37.24"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.85%  '

$ws.Range("D34").Value = "'6.06"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +7.72%  '

$ws.Range("D35").Value = "'53.10"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.83%  '

$ws.Range("D36").Value = "'0.0453"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.30%  '

$ws.Range("E37").Value = '  -0.03%  '

$ws.Range("E38").Value = '  +5.98%  '

$ws.Range("D39").Value = "'18.76"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.68%  '

$ws.Range("E40").Value = '  +1.94%  '

$ws.Range("D41").Value = "'2.71"
$ws.Range("D41").Style = "Normal"

$ws.Range("E42").Value = '  +0.94%  '

$ws.Range("D43").Value = "'22.94"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.79%  '

$ws.Range("E44").Value = '  -1.81%  '

$ws.Range("D45").Value = "'117.34"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.10%  '

$ws.Range("D46").Value = "'3.52"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.61%  '

$ws.Range("D48").Value = '2.186.98'
$ws.Range("E48").Value = '  +3.14%  '

$ws.Range("D49").Value = "'0.251"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +13.56%  '

$ws.Range("D50").Value = "'0.0350"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +10.75%  '

$ws.Range("D51").Value = "'0.953"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.53%  '
